# AttendanceChecking.xlsx edit:
#  - Add a new attendance date column "10/15/2018" by inserting it where the
#    old "Total" header (column N) used to be, and pushing "Total" to the new
#    column O.
#  - Every student data row (6-74): the old per-row COUNTBLANK formula (or
#    stray hard-coded "1") living in column N is cleared, and a fresh
#    COUNTBLANK formula covering E:N (now including the new date column) is
#    written to column O instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 5) -----------------------------------------------
# O5 becomes what N5 used to be ("Total"); copy N5's formatting across too.
$ws.Range("O5").Value = $ws.Range("N5").Value2
$ws.Range("M5").Copy()
$ws.Range("O5").PasteSpecial(-4122)

# N5 becomes the new date column header. Force text (not an auto-converted
# date serial) by formatting the cell as Text before writing the value, then
# restore the original (centered Times New Roman) look by copying the format
# from the still-untouched neighbouring header cell M5.
$ws.Range("N5").NumberFormat = "@"
$ws.Range("N5").Value = "10/15/2018"
$ws.Range("M5").Copy()
$ws.Range("N5").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Data rows (6-74) ---------------------------------------------------
for ($r = 6; $r -le 74; $r++) {
    # Old column-N content (either a =COUNTBLANK(E:H) formula or a stray
    # literal 1) is no longer meaningful now that the Total column moved to
    # O, so just blank it out.
    $ws.Range("N${r}").ClearContents()

    # New Total formula lives in O and counts across the whole E:N span.
    $ws.Range("O${r}").Formula = "=COUNTBLANK(E${r}:N${r})"
    $ws.Range("O${r}").Style = "Normal"
}
